# Figure guide workbook update — "all new variables added 1 YRs"
# Adds a new "Transportation" section (rows 79-83) and a new
# "Public Health and Safety" section (rows 84-94) to the figure guide,
# and clears the (now redundant) yellow highlight on C5:C8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the leftover yellow fill on C5:C8 (Format Cells > Fill > No Color) ---
$ws.Range("C5:C8").Interior.Pattern = -4142

# --- New rows of figure-guide content -------------------------------------
# Values are entered in the same order the original author typed them so the
# shared-string table comes out in the same sequence.

# Row 79 - Transportation / crash data
$ws.Range("H79").Value = "crash data funnel through me from Ashleigh"
$ws.Range("B79").Value = "Transportation"

# Row 80 - Means of Transportation to work
$ws.Range("C80").Value = "Means of Transportation to work"
$ws.Range("H80").Value = "stacked bar for regular geographies all categories"
$ws.Range("B80").Value = "Transportation"

# Row 81 - work from home trend
$ws.Range("H81").Value = "trend line for only work from home last t10 years GOI"

# Row 82 - commute time
$ws.Range("H82").Value = "percent of commutes in GOI over 45 and 90 minutes all 3 years double bar - if 90 isn't significant then drop it or change to a diff category"

# Row 83 - vehicles available
$ws.Range("G83").Value = "vehicles avail"
$ws.Range("H83").Value = "no vehicles, 1 vehicle, 2 or more, leave 3 years and leave county and state, or place county and state - don't need to know incorporated vs unincorporated"

# Row 84 - Public Health and Safety / food insecurity
$ws.Range("B84").Value = "Public Health and Safety"
$ws.Range("H84").Value = "food insecurity same figure type but regular group of geos"

# Row 85 - health insurance
$ws.Range("H85").Value = "w and w/o health insurance figures same types but regular geographies - not incorporated vs unincorporated"

# Row 86 - health insurance by type
$ws.Range("H86").Value = "h insurance by type good fig but regular geographies not incorporated vs unincorporated"

# Row 87 - county health rankings
$ws.Range("H87").Value = "county health rankings - do a gradient table with comparative counties, or just report the GOI"

# Row 88 - premature deaths
$ws.Range("H88").Value = "premature deaths - make something way simpler… add timeline if its available"

# Row 89 - age adjusted death rate
$ws.Range("H89").Value = "age adjusted death rate - just a percent for the GOI for 2010 and current if available "

# Row 90 - leading cause of death
$ws.Range("H90").Value = "leading cause of death regular geographies maybe just top 3, replace the rates - add up all numbers and then find the percent of the total to say ""50% of people dying in x county are dying from 1, 2, 3 at 20%, 20%, and 10%"""

# Row 91 - physical and mental distress
$ws.Range("H91").Value = "physical and mental distress leave as is - put axes on same scale, include a description"

# Row 92 - crime rate per thousand
$ws.Range("H92").Value = "crime rate per thousand replace with regular geographies"

# Row 93 - violent crime rate
$ws.Range("H93").Value = "same ^ with violent crime rate"

# Row 94 - top arrests by type
$ws.Range("H94").Value = "for top arrests by type do one tree plot for GOI"

# --- Leave the cursor where the author left it -----------------------------
$ws.Range("D1").Select() | Out-Null
$ws.Range("A6:A8").Select() | Out-Null
